{"js": "// Update the worksheet date heading and the 24 division-problem answers\n// in the table while leaving the one unchanged cell (65\u00f72=32, 1) intact.\n\nconst body = context.document.body;\n\n// --- 1. Update the date/title paragraph -------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.indexOf(\"2025-08-18 Monday\") !== -1) {\n  titlePara.insertText(\"2025-08-19 Tuesday\", \"Replace\");\n  await context.sync();\n}\n\n// --- 2. Update the table cell contents ---------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (row, col) are 0-indexed. oldText is the pre-edit value (used only to\n// sanity-check we're editing the right cell) and newText is what we write.\nconst updates = [\n  { row: 0, col: 0, oldText: \"46\u00f77=6, 4\", newText: \"98\u00f73=32, 2\" },\n  { row: 0, col: 1, oldText: \"26\u00f75=5, 1\", newText: \"14\u00f76=2, 2\" },\n  { row: 0, col: 2, oldText: \"11\u00f72=5, 1\", newText: \"85\u00f75=17, 0\" },\n  { row: 0, col: 3, oldText: \"87\u00f74=21, 3\", newText: \"55\u00f79=6, 1\" },\n  { row: 0, col: 4, oldText: \"95\u00f77=13, 4\", newText: \"96\u00f74=24, 0\" },\n\n  { row: 4, col: 0, oldText: \"38\u00f74=9, 2\", newText: \"71\u00f74=17, 3\" },\n  { row: 4, col: 1, oldText: \"86\u00f74=21, 2\", newText: \"96\u00f77=13, 5\" },\n  { row: 4, col: 2, oldText: \"88\u00f73=29, 1\", newText: \"13\u00f74=3, 1\" },\n  { row: 4, col: 3, oldText: \"12\u00f74=3, 0\", newText: \"92\u00f74=23, 0\" },\n  { row: 4, col: 4, oldText: \"90\u00f74=22, 2\", newText: \"95\u00f78=11, 7\" },\n\n  { row: 8, col: 0, oldText: \"39\u00f79=4, 3\", newText: \"84\u00f76=14, 0\" },\n  { row: 8, col: 1, oldText: \"15\u00f75=3, 0\", newText: \"73\u00f77=10, 3\" },\n  { row: 8, col: 2, oldText: \"23\u00f79=2, 5\", newText: \"57\u00f73=19, 0\" },\n  { row: 8, col: 3, oldText: \"17\u00f77=2, 3\", newText: \"57\u00f74=14, 1\" },\n  { row: 8, col: 4, oldText: \"35\u00f79=3, 8\", newText: \"68\u00f76=11, 2\" },\n\n  { row: 12, col: 0, oldText: \"59\u00f74=14, 3\", newText: \"82\u00f79=9, 1\" },\n  { row: 12, col: 1, oldText: \"56\u00f76=9, 2\", newText: \"46\u00f72=23, 0\" },\n  { row: 12, col: 2, oldText: \"83\u00f72=41, 1\", newText: \"72\u00f78=9, 0\" },\n  // row 12, col 3 (\"65\u00f72=32, 1\") is intentionally left unchanged.\n  { row: 12, col: 4, oldText: \"59\u00f76=9, 5\", newText: \"75\u00f78=9, 3\" },\n\n  { row: 16, col: 0, oldText: \"86\u00f74=21, 2\", newText: \"26\u00f79=2, 8\" },\n  { row: 16, col: 1, oldText: \"14\u00f78=1, 6\", newText: \"69\u00f72=34, 1\" },\n  { row: 16, col: 2, oldText: \"90\u00f75=18, 0\", newText: \"79\u00f79=8, 7\" },\n  { row: 16, col: 3, oldText: \"67\u00f73=22, 1\", newText: \"73\u00f79=8, 1\" },\n  { row: 16, col: 4, oldText: \"27\u00f75=5, 2\", newText: \"33\u00f77=4, 5\" },\n];\n\n// Verify we are about to overwrite the expected source values (guards\n// against silently writing to the wrong cells if the table shape ever\n// changes), then apply every replacement.\nconst cells = updates.map((u) => table.getCell(u.row, u.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nupdates.forEach((u, i) => {\n  const current = cells[i].value;\n  if (current !== u.oldText) {\n    throw new Error(\n      `Unexpected value at row ${u.row}, col ${u.col}: got ${JSON.stringify(\n        current\n      )}, expected ${JSON.stringify(u.oldText)}`\n    );\n  }\n});\n\nupdates.forEach((u, i) => {\n  cells[i].value = u.newText;\n});\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and the 24 division-problem answers\n# in the table while leaving the one unchanged cell (65\u00f72=32, 1) intact.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date/title paragraph --------------------------------\n$d.Paragraphs.Item(1).Range.Text = \"2025-08-19 Tuesday\"\n\n# --- 2. Update the table cell contents ----------------------------------\n$table = $d.Tables.Item(1)\n\n# (row, col) use 1-based COM indexing; only the five \"content\" rows\n# (1, 5, 9, 13, 17) hold text, the rest are blank spacer rows.\n$table.Cell(1, 1).Range.Text  = \"98\u00f73=32, 2\"\n$table.Cell(1, 2).Range.Text  = \"14\u00f76=2, 2\"\n$table.Cell(1, 3).Range.Text  = \"85\u00f75=17, 0\"\n$table.Cell(1, 4).Range.Text  = \"55\u00f79=6, 1\"\n$table.Cell(1, 5).Range.Text  = \"96\u00f74=24, 0\"\n\n$table.Cell(5, 1).Range.Text  = \"71\u00f74=17, 3\"\n$table.Cell(5, 2).Range.Text  = \"96\u00f77=13, 5\"\n$table.Cell(5, 3).Range.Text  = \"13\u00f74=3, 1\"\n$table.Cell(5, 4).Range.Text  = \"92\u00f74=23, 0\"\n$table.Cell(5, 5).Range.Text  = \"95\u00f78=11, 7\"\n\n$table.Cell(9, 1).Range.Text  = \"84\u00f76=14, 0\"\n$table.Cell(9, 2).Range.Text  = \"73\u00f77=10, 3\"\n$table.Cell(9, 3).Range.Text  = \"57\u00f73=19, 0\"\n$table.Cell(9, 4).Range.Text  = \"57\u00f74=14, 1\"\n$table.Cell(9, 5).Range.Text  = \"68\u00f76=11, 2\"\n\n$table.Cell(13, 1).Range.Text = \"82\u00f79=9, 1\"\n$table.Cell(13, 2).Range.Text = \"46\u00f72=23, 0\"\n$table.Cell(13, 3).Range.Text = \"72\u00f78=9, 0\"\n# $table.Cell(13, 4) (\"65\u00f72=32, 1\") is intentionally left unchanged.\n$table.Cell(13, 5).Range.Text = \"75\u00f78=9, 3\"\n\n$table.Cell(17, 1).Range.Text = \"26\u00f79=2, 8\"\n$table.Cell(17, 2).Range.Text = \"69\u00f72=34, 1\"\n$table.Cell(17, 3).Range.Text = \"79\u00f79=8, 7\"\n$table.Cell(17, 4).Range.Text = \"73\u00f79=8, 1\"\n$table.Cell(17, 5).Range.Text = \"33\u00f77=4, 5\"\n"}
